$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the coin-price / volume(1h) refresh exactly as captured by the diff.
# Column D values are forced to text (leading apostrophe) so Excel does not
# silently reinterpret price strings like "580.16" or "0.0440" as numbers
# (which would drop trailing zeros / change formatting).

$ws.Range("D2").Value = '''69.031.15'
$ws.Range("E2").Value = '  -3.91%  '
$ws.Range("D3").Value = '''3.503.37'
$ws.Range("E3").Value = '  -5.03%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''580.16'
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("D6").Value = '''173.94'
$ws.Range("E6").Value = '  -3.55%  '
$ws.Range("D7").Value = '''0.624'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '''3.497.14'
$ws.Range("E8").Value = '  -4.98%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -6.65%  '
$ws.Range("D11").Value = '''6.69'
$ws.Range("E11").Value = '  +5.93%  '
$ws.Range("D12").Value = '''0.595'
$ws.Range("E12").Value = '  -3.17%  '
$ws.Range("D13").Value = '''47.02'
$ws.Range("E13").Value = '  -6.23%  '
$ws.Range("E14").Value = '  -4.27%  '
$ws.Range("D15").Value = '''674.01'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = '''4.072.57'
$ws.Range("E16").Value = '  -5.03%  '
$ws.Range("D17").Value = '''8.70'
$ws.Range("E17").Value = '  -3.65%  '
$ws.Range("D18").Value = '''69.032.79'
$ws.Range("E18").Value = '  -3.97%  '
$ws.Range("D19").Value = '''3.506.77'
$ws.Range("E19").Value = '  -5.07%  '
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = '''17.48'
$ws.Range("E21").Value = '  -3.74%  '
$ws.Range("D22").Value = '''11.17'
$ws.Range("E22").Value = '  -4.42%  '
$ws.Range("D23").Value = '''0.903'
$ws.Range("E23").Value = '  -4.57%  '
$ws.Range("D24").Value = '''16.09'
$ws.Range("E24").Value = '  -9.92%  '
$ws.Range("D25").Value = '''97.96'
$ws.Range("E25").Value = '  -6.00%  '
$ws.Range("E26").Value = '  -4.77%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  -6.85%  '
$ws.Range("D30").Value = '''9.42'
$ws.Range("E30").Value = '  -7.57%  '
$ws.Range("D31").Value = '''32.78'
$ws.Range("E31").Value = '  -7.79%  '
$ws.Range("D32").Value = '''8.71'
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("D33").Value = '''3.20'
$ws.Range("E33").Value = '  -8.30%  '
$ws.Range("D34").Value = '''7.27'
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("E35").Value = '  -6.28%  '
$ws.Range("D36").Value = '''593.43'
$ws.Range("E36").Value = '  +4.09%  '
$ws.Range("D37").Value = '''3.60'
$ws.Range("E37").Value = '  -15.45%  '
$ws.Range("D38").Value = '''10.89'
$ws.Range("E38").Value = '  -3.91%  '
$ws.Range("E39").Value = '  -4.61%  '
$ws.Range("D40").Value = '''57.26'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '''0.0440'
$ws.Range("E42").Value = '  -5.96%  '
$ws.Range("D43").Value = '''0.336'
$ws.Range("E43").Value = '  -4.68%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '''3.424.73'
$ws.Range("E44").Value = '  -9.79%  '
$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D45").Value = '''0.136'
$ws.Range("E45").Value = '  -7.45%  '
$ws.Range("D46").Value = '''33.39'
$ws.Range("E46").Value = '  -6.18%  '
$ws.Range("D47").Value = '''0.0₃0705'
$ws.Range("E47").Value = '  -9.93%  '
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").Value = '''2.59'
$ws.Range("E49").Value = '  -8.08%  '
$ws.Range("D50").Value = '''0.133'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").Value = '''5.76'
$ws.Range("E51").Value = '  +17.84%  '
